$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nuevas ciudades del AMB (Área Metropolitana de Bucaramanga) con población
# proyectada DANE 2020.
$ws.Cells.Item(12, 1).Value = "FLORIDABLANCA"
$ws.Cells.Item(12, 2).Value = 307896

$ws.Cells.Item(13, 1).Value = "GIRON"
$ws.Cells.Item(13, 2).Value = 171904

$ws.Cells.Item(14, 1).Value = "PIEDECUESTA"
$ws.Cells.Item(14, 2).Value = 182959

# Replicar el formato numérico/estilo de la columna B a las filas nuevas.
$ws.Range("B11").Copy()
$ws.Range("B12:B14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A16").Select()
